$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in the order the per-row values appear): B, C, D, E, F, I, J, K, L, N
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 14)

$data = @(
    ,@(2, 1.545629258784857, 0.07717628820958566, 0.01815585547011267, 0.06961748551169222, 3.570975794395295, 2.218182416072231, 0.1477815485477514, 1.352588289995509, 0.3814899393343723, 3.514942601426007)
    ,@(3, 1.51096712221559, 0.07145143369969276, 0.01780671062417127, 0.06968500031629077, 3.570356853167937, 2.221860651120018, 0.148324769332369, 1.314252251397221, 0.3783461650575219, 3.534378708461389)
    ,@(4, 1.490490830023106, 0.06797973033009441, 0.01759941358067252, 0.06974613591110135, 3.571361696658784, 2.224912492693676, 0.1486949130674198, 1.291455622419477, 0.3765853697630277, 3.547150963245514)
    ,@(5, 1.482349752132905, 0.06657583476550144, 0.01751673526084829, 0.06977601291381141, 3.57211958177399, 2.226355728356701, 0.148854970079709, 1.282352599688409, 0.375910529532355, 3.552566584178678)
    ,@(6, 1.481010215716026, 0.06634337275414737, 0.01750311561294993, 0.06978127420067359, 3.57226647569162, 2.226607433454376, 0.148882104790367, 1.280852338585163, 0.3758010539815544, 3.55347857716626)
    ,@(7, 1.490380213408997, 0.06796075303942928, 0.01759829125448675, 0.06974651872506143, 3.571370506881053, 2.224931148459788, 0.1486970342977543, 1.291332099247143, 0.3765760956449284, 3.547223146601162)
    ,@(8, 1.533510610706855, 0.07519331871972668, 0.01803401087222767, 0.06963668851190086, 3.570475042491807, 2.219286048031428, 0.1479612634908918, 1.339216181493526, 0.3803708316214198, 3.521470159716721)
    ,@(9, 1.624476395604432, 0.08972391575447602, 0.01894395548658423, 0.06957686567079158, 3.579703873641066, 2.214509288776469, 0.1468081733529871, 1.438999072899236, 0.3891545921644877, 3.477620897697577)
    ,@(10, 1.695197219962267, 0.1006172201936408, 0.01964551899122924, 0.06962698538396772, 3.593181707804931, 2.214836170887708, 0.1461367586338511, 1.515900363879126, 0.3964240867530151, 3.449459894758121)
    ,@(11, 1.728213758628101, 0.1056214530945283, 0.01997168734945021, 0.0696700574089224, 3.600768040855087, 2.215818099032859, 0.1458692977201004, 1.551666185535822, 0.3999079873146769, 3.437529257273937)
    ,@(12, 1.740837575793989, 0.1075235272834618, 0.02009619368610061, 0.06968926908403894, 3.603850007860089, 2.216309729254178, 0.1457734622182336, 1.565322305470545, 0.4012526427625289, 3.433138004589608)
    ,@(13, 1.738113427806695, 0.1071135657327602, 0.02006933509577991, 0.0696850026932605, 3.603176948491708, 2.216198520272599, 0.1457938601069237, 1.562376221061157, 0.4009619195746836, 3.434078106512729)
    ,@(14, 1.729249901385572, 0.1057777956453378, 0.01998191074508071, 0.06967157989918427, 3.601017404084274, 2.215856144984585, 0.1458613042048462, 1.552787431737528, 0.4000181047159401, 3.43716544809665)
    ,@(15, 1.72383650430487, 0.1049605216954888, 0.01992848969486616, 0.06966373546561933, 3.599721861198589, 2.215662030527852, 0.1459033245048325, 1.546928650059897, 0.3994432937436017, 3.439073026349632)
    ,@(16, 1.693056488154753, 0.1002911692206112, 0.01962434303742455, 0.06962457726610438, 3.592715199444072, 2.214788764315415, 0.1461550006880827, 1.513578721115181, 0.396199960740617, 3.450257303944824)
    ,@(17, 1.674390170281072, 0.0974392180087591, 0.01943954600863051, 0.06960573728445851, 3.588789500354764, 2.214466436963654, 0.1463191108168154, 1.493320016686312, 0.3942555572566846, 3.457343930555353)
    ,@(18, 1.66373338221365, 0.09580344208092129, 0.01933391771606452, 0.06959680997720241, 3.586668515986759, 2.214359471230537, 0.1464170772580395, 1.481741461946683, 0.3931538493068274, 3.461502789782699)
    ,@(19, 1.660138857207073, 0.09525038312901302, 0.01929826800918377, 0.06959411563359197, 3.58597391471163, 2.214336725354102, 0.1464508613809912, 1.47783382988132, 0.3927836937083953, 3.462925129776238)
    ,@(20, 1.676368997775967, 0.09774233750007966, 0.01945914956678152, 0.06960754536447489, 3.589193221939297, 2.214492632371297, 0.1463012711517919, 1.49546896261424, 0.3944608184354195, 3.456580974971942)
    ,@(21, 1.731850048422757, 0.1061699510729568, 0.02000756256414249, 0.06967544387091529, 3.601646037860377, 2.215953457853175, 0.1458413465409052, 1.555600843386912, 0.4002946378256809, 3.436255183875971)
    ,@(22, 1.768816223718261, 0.1117191785284319, 0.02037176525657713, 0.06973672370039985, 3.611003948694119, 2.217606442728737, 0.1455724972210781, 1.595555437929931, 0.4042552658728482, 3.423709190280888)
    ,@(23, 1.749022196283477, 0.1087536528538635, 0.0201768596844154, 0.06970247533132756, 3.605897913040693, 2.21666033001263, 0.1457130876410879, 1.574171064109152, 0.4021278952440355, 3.430337658783785)
    ,@(24, 1.675474137066317, 0.09760528519461786, 0.01945028488424683, 0.06960672199987705, 3.589010275904172, 2.214480545373917, 0.1463093251954177, 1.494497210945013, 0.3943679695293412, 3.456925643451115)
    ,@(25, 1.599184536871718, 0.08575519842256085, 0.01869193384603562, 0.06957647496785313, 3.5760311691688, 2.215127821023948, 0.1470891862454273, 1.411374934514924, 0.3866348840376048, 3.488771166570814)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $c = $cols[$i]
        $val = $rowData[$i + 1]
        $ws.Cells.Item($r, $c).Value = $val
    }
}

Write-Output "done"
